$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 (header/ID) values updated
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) values updated
$ws.Range("B2").Value = 15.343105725750167
$ws.Range("C2").Value = 10.65406927831711
$ws.Range("D2").Value = 15.380396980945317
$ws.Range("E2").Value = 7.8425795026476139

# Row 3 (STR) values updated
$ws.Range("B3").Value = 13.133227646912529
$ws.Range("C3").Value = 12.959001009682119
$ws.Range("D3").Value = 11.131336492850405
$ws.Range("E3").Value = 13.473689143994417

# Update selection to reflect the edited range
$ws.Range("B1:E3").Select()
